# Update training_history values for run_7 (128 dense layers, 20 lstm, 50 epochs, 0.5 dropout)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 50,4
$data[0,0] = 0.06430843472480774
$data[0,1] = 0.9817928075790405
$data[0,2] = 0.1682728081941605
$data[0,3] = 0.9220559597015381
$data[1,0] = 0.01029054820537567
$data[1,1] = 0.998230516910553
$data[1,2] = 0.07607636600732803
$data[1,3] = 0.9603598713874817
$data[2,0] = 0.005895404145121574
$data[2,1] = 0.9984866380691528
$data[2,2] = 0.06404083967208862
$data[2,3] = 0.9717620015144348
$data[3,0] = 0.002480259397998452
$data[3,1] = 0.999371349811554
$data[3,2] = 0.03020020946860313
$data[3,3] = 0.9926955103874207
$data[4,0] = 0.001659055822528899
$data[4,1] = 0.9996041655540466
$data[4,2] = 0.036944430321455
$data[4,3] = 0.9891323447227478
$data[5,0] = 0.001108076306991279
$data[5,1] = 0.9997206330299377
$data[5,2] = 0.006993710063397884
$data[5,3] = 0.9992873668670654
$data[6,0] = 0.001555933384224772
$data[6,1] = 0.9996274709701538
$data[6,2] = 0.07404659688472748
$data[6,3] = 0.9722073674201965
$data[7,0] = 0.001035849447362125
$data[7,1] = 0.9997438788414001
$data[7,2] = 0.04312290251255035
$data[7,3] = 0.9770176410675049
$data[8,0] = 0.001184522756375372
$data[8,1] = 0.9997206330299377
$data[8,2] = 0.01034195721149445
$data[8,3] = 0.9982184171676636
$data[9,0] = 0.000480839895317331
$data[9,1] = 0.9998602867126465
$data[9,2] = 0.04988864064216614
$data[9,3] = 0.9836985468864441
$data[10,0] = 0.0008925613365136087
$data[10,1] = 0.9998370409011841
$data[10,2] = 0.1782244443893433
$data[10,3] = 0.9297167062759399
$data[11,0] = 0.000338803103659302
$data[11,1] = 0.9999068975448608
$data[11,2] = 0.2251095473766327
$data[11,3] = 0.933368980884552
$data[12,0] = 0.000606511312071234
$data[12,1] = 0.9997904300689697
$data[12,2] = 0.08357534557580948
$data[12,3] = 0.9754142165184021
$data[13,0] = 0.0005029302556067705
$data[13,1] = 0.9998370409011841
$data[13,2] = 0.04891286417841911
$data[13,3] = 0.9712275266647339
$data[14,0] = 0.0003233772586099803
$data[14,1] = 0.9999068975448608
$data[14,2] = 0.09727617353200912
$data[14,3] = 0.9589346051216125
$data[15,0] = 0.0001811326656024903
$data[15,1] = 0.9999534487724304
$data[15,2] = 0.5239173173904419
$data[15,3] = 0.9210760593414307
$data[16,0] = 0.001151273492723703
$data[16,1] = 0.9997206330299377
$data[16,2] = 0.01319500431418419
$data[16,3] = 0.9947443604469299
$data[17,0] = 0.000031709252652945
$data[17,1] = 1
$data[17,2] = 0.02498102746903896
$data[17,3] = 0.9927846193313599
$data[18,0] = 0.0006060707964934409
$data[18,1] = 0.9998602867126465
$data[18,2] = 0.04773042351007462
$data[18,3] = 0.9916265606880188
$data[19,0] = 0.0003361108538229018
$data[19,1] = 0.9999301433563232
$data[19,2] = 0.02034757845103741
$data[19,3] = 0.9937644600868225
$data[20,0] = 0.0002420740784145892
$data[20,1] = 0.9999068975448608
$data[20,2] = 0.01323843933641911
$data[20,3] = 0.9981293678283691
$data[21,0] = 0.0004544209514278919
$data[21,1] = 0.9999068975448608
$data[21,2] = 0.01054247654974461
$data[21,3] = 0.9967931509017944
$data[22,0] = 0.0002120501449098811
$data[22,1] = 0.9999534487724304
$data[22,2] = 0.04884923249483109
$data[22,3] = 0.9755033254623413
$data[23,0] = 0.0002107933833030984
$data[23,1] = 0.9999068975448608
$data[23,2] = 0.03461636230349541
$data[23,3] = 0.9844111800193787
$data[24,0] = 0.0001322500174865127
$data[24,1] = 0.9999534487724304
$data[24,2] = 0.003626617603003979
$data[24,3] = 0.9984856843948364
$data[25,0] = 0.0004938208730891347
$data[25,1] = 0.9998602867126465
$data[25,2] = 0.1122124865651131
$data[25,3] = 0.9405843615531921
$data[26,0] = 0.0004936656914651394
$data[26,1] = 0.9999301433563232
$data[26,2] = 0.03752052411437035
$data[26,3] = 0.9877961874008179
$data[27,0] = 0.00005018249794375151
$data[27,1] = 0.9999766945838928
$data[27,2] = 0.1134463474154472
$data[27,3] = 0.9443256855010986
$data[28,0] = 0.0004552112077362835
$data[28,1] = 0.9999068975448608
$data[28,2] = 0.03929530456662178
$data[28,3] = 0.9926955103874207
$data[29,0] = 0.0001059204005287029
$data[29,1] = 0.9999534487724304
$data[29,2] = 0.0772145539522171
$data[29,3] = 0.9837876558303833
$data[30,0] = 0.00005000036617275327
$data[30,1] = 1
$data[30,2] = 0.02384382672607899
$data[30,3] = 0.9941207766532898
$data[31,0] = 0.000006125530489953235
$data[31,1] = 1
$data[31,2] = 0.04749710485339165
$data[31,3] = 0.9927846193313599
$data[32,0] = 0.00003625395402195863
$data[32,1] = 0.9999766945838928
$data[32,2] = 0.03064015135169029
$data[32,3] = 0.9920719861984253
$data[33,0] = 0.0002965771127492189
$data[33,1] = 0.9999766945838928
$data[33,2] = 0.001772635034285486
$data[33,3] = 0.9993764758110046
$data[34,0] = 0.0004728272324427962
$data[34,1] = 0.9998835921287537
$data[34,2] = 0.001001277123577893
$data[34,3] = 0.9998218417167664
$data[35,0] = 0.00001894428169180173
$data[35,1] = 1
$data[35,2] = 0.0007723688031546772
$data[35,3] = 0.9999109506607056
$data[36,0] = 0.00013053446309641
$data[36,1] = 0.9999534487724304
$data[36,2] = 0.1395634263753891
$data[36,3] = 0.9566185474395752
$data[37,0] = 0.0003434133250266314
$data[37,1] = 0.9999068975448608
$data[37,2] = 0.004451545886695385
$data[37,3] = 0.9994655251502991
$data[38,0] = 0.0006022296147421002
$data[38,1] = 0.9998602867126465
$data[38,2] = 0.08579479902982712
$data[38,3] = 0.9837876558303833
$data[39,0] = 0.0001212255810969509
$data[39,1] = 0.9999534487724304
$data[39,2] = 0.466886967420578
$data[39,3] = 0.9306966066360474
$data[40,0] = 0.0001790060778148472
$data[40,1] = 0.9999534487724304
$data[40,2] = 0.3190890252590179
$data[40,3] = 0.9259754419326782
$data[41,0] = 0.00000245676619670121
$data[41,1] = 1
$data[41,2] = 0.1105284839868546
$data[41,3] = 0.9461072683334351
$data[42,0] = 0.000001847836415436177
$data[42,1] = 1
$data[42,2] = 0.1132524237036705
$data[42,3] = 0.947443425655365
$data[43,0] = 0.0004547676362562925
$data[43,1] = 0.9998137354850769
$data[43,2] = 0.01490586157888174
$data[43,3] = 0.9962586760520935
$data[44,0] = 0.0001494372118031606
$data[44,1] = 0.9999534487724304
$data[44,2] = 0.08132679760456085
$data[44,3] = 0.9711384177207947
$data[45,0] = 0.0004402854247018695
$data[45,1] = 0.9999301433563232
$data[45,2] = 0.01216733176261187
$data[45,3] = 0.9958133101463318
$data[46,0] = 0.0003573191352188587
$data[46,1] = 0.9998835921287537
$data[46,2] = 0.005245896056294441
$data[46,3] = 0.9994655251502991
$data[47,0] = 0.00009431406215298921
$data[47,1] = 0.9999766945838928
$data[47,2] = 0.1915923058986664
$data[47,3] = 0.9362195134162903
$data[48,0] = 0.000002520067027944606
$data[48,1] = 1
$data[48,2] = 0.2152108550071716
$data[48,3] = 0.9323000311851501
$data[49,0] = 0.000002108103444697917
$data[49,1] = 1
$data[49,2] = 0.2482743561267853
$data[49,3] = 0.9317655563354492

$ws.Range("A2:D51").Value = $data
